# "ReUpload new dates for italian relay"
#
# 2020-11-15 originally had a single "OT" (Team Open) entry split by
# gender (SplitGender=1) with an Evening (S) and Final (F) heat at
# 11:30/13:00. The re-upload replaces that with separate "U17T" and
# "U13T" team events, not split by gender, all at 10:00 (rollcall
# 10:45): U17T keeps the existing row 105/106 slots (E + F heats) and
# U13T is appended as two brand-new rows (107/108, E + F heats). Row
# 104 (the surviving "OT" row) also flips to SplitGender=0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

$kFormula = '=IF(Table1[[#This Row],[Cancelled]]=1,"N/A",Table1[[#This Row],[Date]]&Table1[[#This Row],[Category]]&IF(Table1[[#This Row],[SplitGender]]=1,IF(OR(Table1[[#This Row],[Category]]="U9",Table1[[#This Row],[Category]]="U11",Table1[[#This Row],[Category]]="U13"),"B","M"),"")&Table1[[#This Row],[Weapon]])'
$lFormula = '=IF(Table1[[#This Row],[Cancelled]]=1,"",IF(Table1[[#This Row],[SplitGender]]=0,"N/A",Table1[[#This Row],[Date]]&Table1[[#This Row],[Category]]&IF(Table1[[#This Row],[SplitGender]]=1,IF(OR(Table1[[#This Row],[Category]]="U9",Table1[[#This Row],[Category]]="U11",Table1[[#This Row],[Category]]="U13"),"G","W"),"")&Table1[[#This Row],[Weapon]]))'

# --- Row 104: 20201115 OT / E  -> SplitGender 1 -> 0 -------------------
$ws.Cells(104, 5).Value2 = 0

# --- Row 105: 20201115 OT / S  -> U17T / F, SplitGender 0, new time ----
$ws.Cells(105, 2).Value2 = "U17T"
$ws.Cells(105, 3).Value2 = "F"
$ws.Cells(105, 5).Value2 = 0
$ws.Cells(105, 6).Value2 = "10:00"
$ws.Cells(105, 7).Value2 = "10:45"

# --- Row 106: 20201115 OT / F  -> U17T / E, SplitGender 0, new time ----
$ws.Cells(106, 2).Value2 = "U17T"
$ws.Cells(106, 3).Value2 = "E"
$ws.Cells(106, 5).Value2 = 0
$ws.Cells(106, 6).Value2 = "10:00"
$ws.Cells(106, 7).Value2 = "10:45"

# --- New row: 20201115 U13T / F -----------------------------------------
$newRow1 = $lo.ListRows.Add()
$r1 = $newRow1.Range.Row
$ws.Cells($r1, 1).Value2 = 20201115
$ws.Cells($r1, 2).Value2 = "U13T"
$ws.Cells($r1, 2).NumberFormat = "@"
$ws.Cells($r1, 3).Value2 = "F"
$ws.Cells($r1, 3).NumberFormat = "@"
$ws.Cells($r1, 4).Value2 = 0
$ws.Cells($r1, 5).Value2 = 0
$ws.Cells($r1, 6).Value2 = "10:00"
$ws.Cells($r1, 6).NumberFormat = "@"
$ws.Cells($r1, 7).Value2 = "10:45"
$ws.Cells($r1, 7).NumberFormat = "@"
$ws.Cells($r1, 8).Value2 = "FSA"
$ws.Cells($r1, 8).NumberFormat = "@"
$ws.Cells($r1, 9).NumberFormat = "@"
$ws.Cells($r1, 11).Formula = $kFormula
$ws.Cells($r1, 12).Formula = $lFormula

# --- New row: 20201115 U13T / E -----------------------------------------
$newRow2 = $lo.ListRows.Add()
$r2 = $newRow2.Range.Row
$ws.Cells($r2, 1).Value2 = 20201115
$ws.Cells($r2, 2).Value2 = "U13T"
$ws.Cells($r2, 2).NumberFormat = "@"
$ws.Cells($r2, 3).Value2 = "E"
$ws.Cells($r2, 3).NumberFormat = "@"
$ws.Cells($r2, 4).Value2 = 0
$ws.Cells($r2, 5).Value2 = 0
$ws.Cells($r2, 6).Value2 = "10:00"
$ws.Cells($r2, 6).NumberFormat = "@"
$ws.Cells($r2, 7).Value2 = "10:45"
$ws.Cells($r2, 7).NumberFormat = "@"
$ws.Cells($r2, 8).Value2 = "FSA"
$ws.Cells($r2, 8).NumberFormat = "@"
$ws.Cells($r2, 9).NumberFormat = "@"
$ws.Cells($r2, 11).Formula = $kFormula
$ws.Cells($r2, 12).Formula = $lFormula
